# Commit: "Adding 3D Steel to the portfolio forecast"
# The data window rolled forward by 14 days (Entsoe actual-wind-production
# pull refreshed), and column B ("Actual Production (MW)") was replaced
# with the new period's readings.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds Excel date-time serials; every row's timestamp shifts
# forward by exactly 14 days (45875.x -> 45889.x, ..., 45876 -> 45890).
for ($r = 2; $r -le 97; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $cur + 14
}

# Column B holds the new Actual Production (MW) readings for rows 2..97.
$newB = @(
    226,262,301,370,465,523,548,589,582,575,583,583,568,546,527,488,
    440,424,399,371,333,307,294,287,256,238,232,231,202,141,103,75,
    51,35,24,19,25,27,32,44,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,
    0,0
)

for ($i = 0; $i -lt $newB.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 2).Value = $newB[$i]
}
